# Update countries & provincias Spain
# Refreshes the COVID "Pais" dataset snapshot: a handful of countries
# got revised totals (and, as a side effect of the new totals, a few
# adjacent rows swap rank/name because the sheet is kept sorted by
# "Casos totales" descending), plus the "last updated" timestamp moves
# from 10:32 to 11:49.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: "last updated" banner text
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 29 de Julio de 2020 a las 11:49"

# Row 4 - Estados Unidos (name/B/C unchanged)
$ws.Cells.Item(4, 4).Value = 2188415   # Casos activos (D)
$ws.Cells.Item(4, 5).Value = 2157587   # Recuperados (E)
$ws.Cells.Item(4, 7).Value = 21        # Muertes hoy (G)
$ws.Cells.Item(4, 8).Value = 152341    # Muertes (H)

# Row 19 - Banglades
$ws.Cells.Item(19, 2).Value = 232194   # Casos totales (B)
$ws.Cells.Item(19, 3).Value = 3009     # Nuevos casos (C)
$ws.Cells.Item(19, 4).Value = 130292   # Casos activos (D)
$ws.Cells.Item(19, 5).Value = 98867    # Recuperados (E)
$ws.Cells.Item(19, 7).Value = 35       # Muertes hoy (G)
$ws.Cells.Item(19, 8).Value = 3035     # Muertes (H)

# Row 30 - now Filipinas (was China); Filipinas overtakes China in rank
$ws.Cells.Item(30, 1).Value = "Filipinas"
$ws.Cells.Item(30, 2).Value = 85486
$ws.Cells.Item(30, 3).Value = 1874
$ws.Cells.Item(30, 4).Value = 26996
$ws.Cells.Item(30, 5).Value = 56528
$ws.Cells.Item(30, 7).Value = 16
$ws.Cells.Item(30, 8).Value = 1962

# Row 31 - now China (was Filipinas), takes the previous China figures
$ws.Cells.Item(31, 1).Value = "China"
$ws.Cells.Item(31, 2).Value = 84060
$ws.Cells.Item(31, 3).Value = 101
$ws.Cells.Item(31, 4).Value = 78944
$ws.Cells.Item(31, 5).Value = 482
$ws.Cells.Item(31, 8).Value = 4634

# Row 34 - Oman
$ws.Cells.Item(34, 2).Value = 78569
$ws.Cells.Item(34, 3).Value = 665
$ws.Cells.Item(34, 4).Value = 60240
$ws.Cells.Item(34, 5).Value = 17917
$ws.Cells.Item(34, 7).Value = 10
$ws.Cells.Item(34, 8).Value = 412

# Row 38 - now Israel (was Belgica); Israel overtakes Belgica in rank
$ws.Cells.Item(38, 1).Value = "Israel"
$ws.Cells.Item(38, 2).Value = 66805
$ws.Cells.Item(38, 3).Value = 512
$ws.Cells.Item(38, 4).Value = 32697
$ws.Cells.Item(38, 5).Value = 33618
$ws.Cells.Item(38, 7).Value = 4
$ws.Cells.Item(38, 8).Value = 490

# Row 39 - now Belgica (was Israel), takes the previous Belgica figures
$ws.Cells.Item(39, 1).Value = "Belgica"
$ws.Cells.Item(39, 2).Value = 66662
$ws.Cells.Item(39, 3).Value = 234
$ws.Cells.Item(39, 4).Value = 17476
$ws.Cells.Item(39, 5).Value = 39353
$ws.Cells.Item(39, 7).Value = 11
$ws.Cells.Item(39, 8).Value = 9833

# Row 49 - Polonia
$ws.Cells.Item(49, 2).Value = 44416
$ws.Cells.Item(49, 3).Value = 512
$ws.Cells.Item(49, 4).Value = 33190
$ws.Cells.Item(49, 5).Value = 9532
$ws.Cells.Item(49, 7).Value = 12
$ws.Cells.Item(49, 8).Value = 1694

# Row 66 - Austria
$ws.Cells.Item(66, 2).Value = 20850
$ws.Cells.Item(66, 3).Value = 173
$ws.Cells.Item(66, 4).Value = 18528
$ws.Cells.Item(66, 5).Value = 1606
$ws.Cells.Item(66, 7).Value = 3
$ws.Cells.Item(66, 8).Value = 716

# Row 87 - Malasia
$ws.Cells.Item(87, 2).Value = 8956
$ws.Cells.Item(87, 3).Value = 13
$ws.Cells.Item(87, 4).Value = 8612
$ws.Cells.Item(87, 5).Value = 220

# Row 88 - Consejo Danes para los Refugiados
$ws.Cells.Item(88, 2).Value = 8931
$ws.Cells.Item(88, 3).Value = 58
$ws.Cells.Item(88, 4).Value = 6095
$ws.Cells.Item(88, 5).Value = 2626
$ws.Cells.Item(88, 7).Value = 2
$ws.Cells.Item(88, 8).Value = 210

# Row 90 - Finlandia
$ws.Cells.Item(90, 2).Value = 7414
$ws.Cells.Item(90, 3).Value = 10
$ws.Cells.Item(90, 5).Value = 165

# Row 114 - now Hong Kong (was Montenegro); Hong Kong jumps to the top of this block
$ws.Cells.Item(114, 1).Value = "Hong Kong"
$ws.Cells.Item(114, 2).Value = 3003
$ws.Cells.Item(114, 3).Value = 118
$ws.Cells.Item(114, 4).Value = 1527
$ws.Cells.Item(114, 5).Value = 1452
$ws.Cells.Item(114, 7).Value = 1
$ws.Cells.Item(114, 8).Value = 24

# Row 115 - now Montenegro (was Mayotte), takes the previous Montenegro figures
$ws.Cells.Item(115, 1).Value = "Montenegro"
$ws.Cells.Item(115, 2).Value = 2949
$ws.Cells.Item(115, 4).Value = 839
$ws.Cells.Item(115, 5).Value = 2065
$ws.Cells.Item(115, 8).Value = 45

# Row 116 - now Mayotte (was Hong Kong), takes the previous Mayotte figures
$ws.Cells.Item(116, 1).Value = "Mayotte"
$ws.Cells.Item(116, 2).Value = 2900
$ws.Cells.Item(116, 4).Value = 2672
$ws.Cells.Item(116, 5).Value = 190
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 38

# Row 118 - Sri Lanka
$ws.Cells.Item(118, 4).Value = 2317
$ws.Cells.Item(118, 5).Value = 482

# Row 125 - Eslovenia
$ws.Cells.Item(125, 2).Value = 2115
$ws.Cells.Item(125, 3).Value = 14
$ws.Cells.Item(125, 4).Value = 1761
$ws.Cells.Item(125, 5).Value = 237
